$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" footer field text (2022/10/8 ->
#    2022/12/31) on the slide master and every slide layout. PowerPoint
#    re-caches this text whenever the deck is saved on a later date.
# ---------------------------------------------------------------------------
$newDate = "2022/12/31"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 8 ("...你真偉大"): split the lyric run
#    "當主再來  歡呼聲響徹天空" into three runs, changing "響徹天" to
#    "響澈天" in the process: "當主再來  歡呼聲" + "響澈天" + "空".
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
for ($k = 1; $k -le $slide8.Shapes.Count; $k++) {
    $shape = $slide8.Shapes.Item($k)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text.IndexOf("當主再來") -ge 0) {
            # Peel the trailing "空" off into its own run first (content
            # unchanged) so the 3-character middle chunk below splits
            # cleanly into its own run too.
            $tail = $tr.Characters(13, 1)
            $tail.Text = "空"

            # Re-text the middle 3 characters ("響徹天" -> "響澈天") as a
            # separate run.
            $mid = $tr.Characters(10, 3)
            $mid.Text = "響澈天"
        }
    }
}
